$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- New row 8: second "New Storage Processes" technology definition ---
# Copy formatting from row 7 (the existing NEW_ELEC_STG template row)
$ws.Range("B7:P7").Copy()
$ws.Range("B8:P8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B8").Value = "NEW_ELEC_STG2"
$ws.Range("C8").Value = "new storage"
$ws.Range("E8").Formula = "=E7"
$ws.Range("G8").Formula = "=G7"
$ws.Range("H8").Value = 2027
$ws.Range("I8").Value = 0.8
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 15
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3500
$ws.Range("N8").Value = 0.8
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 1

# --- New row 22: second "STS" commodity-set membership entry ---
# Copy formatting from row 21 (the existing STS/NEW_ELEC_STG membership row)
$ws.Range("B21:I21").Copy()
$ws.Range("B22:I22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B22").Value = "STG"
$ws.Range("C22").Formula = "=B8"
$ws.Range("D22").Formula = "=C8"
$ws.Range("E22").Formula = "=E21"
$ws.Range("F22").Formula = "=F21"
$ws.Range("G22").Formula = "=G21"
$ws.Range("H22").Value = $null
$ws.Range("I22").Value = $null

# --- Cosmetic view-state refresh (zoom + selection), as left by the editing session ---
$ws.Range("H25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 83
